$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Move the bottom-most rows out of the way first (process bottom-up
#    so we never overwrite data we still need to read).
# ---------------------------------------------------------------------

# row 34 (B34 "end") -> row 36
$ws.Range("B36").Value = "end"
$ws.Rows.Item(36).RowHeight = 12.8
$ws.Range("B34").Clear()

# row 33 (H/I/J/K, step 10) -> row 34
$ws.Range("H34").Value = 10
$ws.Range("I33").Copy()
$ws.Range("I34").PasteSpecial(-4122)
$ws.Range("I34").Value = "あはは、天使なんかいないの。でも、大丈夫、#bigdaddyがいれば。"
$ws.Range("J33").Copy()
$ws.Range("J34").PasteSpecial(-4122)
$ws.Range("J34").Value = "Ahaha, there are no angels. But it's okay, as long as #bigdaddy is here."
$ws.Range("K34").Value = "啊哈哈，天使都是假的。但没关系，只要有#bigdaddy在就好。"
$ws.Rows.Item(34).RowHeight = 23.85
$ws.Range("H33:K33").Clear()

# row 32 (H/I/J/K, step 9) -> row 33
$ws.Range("H33").Value = 9
$ws.Range("I32").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("I33").Value = "…リトル…私たちのリトル…。天使は残酷ね、#bigdaddy…天使は私たちのリトルを見殺しにしたのよ。私たちが…リトルを…この広い世界から守らなきゃ。彼女たちのための…リトル…ガーデン…"
$ws.Range("J32").Copy()
$ws.Range("J33").PasteSpecial(-4122)
$ws.Range("J33").Value = "...Our Little... our Little One... Angels are cruel, #bigdaddy...Angels let our Little die. We must...protect them...from this vast world. Here..for them... the Little..Garden..."
$ws.Range("K33").Value = "…小妹妹…我们的小妹妹…。天使真是残酷，#bigdaddy…天使对小妹妹的死冷眼旁观。我们…必须从这广阔的世界中…保护小妹妹。这里是为她们而存在的…小小…花园…"
$ws.Rows.Item(33).RowHeight = 57.45
$ws.Range("H32:K32").Clear()

# row 31 (A31 "little_dead") -> row 32
$ws.Range("A32").Value = "little_dead"
$ws.Rows.Item(32).RowHeight = 12.8
$ws.Range("A31").Clear()

# row 28 (B28 "end") -> row 29
$ws.Range("B29").Value = "end"
$ws.Rows.Item(29).RowHeight = 12.8
$ws.Range("B28").Clear()

# ---------------------------------------------------------------------
# 2) Insert the two new "modAffinity" rows (27 and 35).
# ---------------------------------------------------------------------

# row 27: modAffinity, 100 (right after row 26, same row height 91)
$ws.Range("D27").Value = "modAffinity"
$ws.Range("E27").Value = 100
$ws.Range("I26").Copy()
$ws.Range("I27").PasteSpecial(-4122)
$ws.Range("J26").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Rows.Item(27).RowHeight = 91

# row 35: modAffinity, -200 (between the relocated rows 34 and 36)
$ws.Range("D35").Value = "modAffinity"
$ws.Range("E35").Value = -200
$ws.Range("I34").Copy()
$ws.Range("I35").PasteSpecial(-4122)
$ws.Range("J34").Copy()
$ws.Range("J35").PasteSpecial(-4122)
$ws.Rows.Item(35).RowHeight = 13.8

# ---------------------------------------------------------------------
# 3) Sheet view bookkeeping (pane/selection) to mirror the saved state.
# ---------------------------------------------------------------------
$ws.Application.Goto($ws.Range("A23"))
$ws.Range("K34").Select()
$ws.Range("H37").Select()

Write-Output "edit complete"
